# Remove duplicate entry "Spectra.jl" (row 135) from the FOSS4Spec sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire duplicate row and delete it, shifting the rows below up.
$row = $ws.Rows.Item(135)
$row.Select()
$row.EntireRow.Delete()
